$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = "2025-07-18T08:10:51Z"
$ws.Range("R2").Value = 8441
$ws.Range("S2").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5303487', 'current_grade': None, 'current_score': 86.15, 'final_grade': None, 'final_score': 73.05, 'unposted_current_score': 86.15, 'unposted_current_grade': None, 'unposted_final_score': 73.05, 'unposted_final_grade': None}"

$ws.Range("P3").Value = "2025-07-17T15:19:05Z"
$ws.Range("R3").Value = 24153
$ws.Range("S3").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5029701', 'current_grade': None, 'current_score': 98.77, 'final_grade': None, 'final_score': 88.9, 'unposted_current_score': 98.77, 'unposted_current_grade': None, 'unposted_final_score': 88.9, 'unposted_final_grade': None}"

$ws.Range("P4").Value = "2025-07-21T00:34:17Z"
$ws.Range("R4").Value = 23313
$ws.Range("S4").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5211091', 'current_grade': None, 'current_score': 99.36, 'final_grade': None, 'final_score': 89.42, 'unposted_current_score': 99.36, 'unposted_current_grade': None, 'unposted_final_score': 89.42, 'unposted_final_grade': None}"

$ws.Range("P5").Value = "2025-07-21T01:12:44Z"
$ws.Range("R5").Value = 12513
$ws.Range("S5").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5252563', 'current_grade': None, 'current_score': 98.27, 'final_grade': None, 'final_score': 81.77, 'unposted_current_score': 98.27, 'unposted_current_grade': None, 'unposted_final_score': 81.77, 'unposted_final_grade': None}"

$ws.Range("P6").Value = "2025-07-20T01:53:47Z"
$ws.Range("R6").Value = 21940
$ws.Range("S6").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5243067', 'current_grade': None, 'current_score': 96.81, 'final_grade': None, 'final_score': 87.13, 'unposted_current_score': 96.81, 'unposted_current_grade': None, 'unposted_final_score': 87.13, 'unposted_final_grade': None}"

$ws.Range("P7").Value = "2025-07-15T21:05:26Z"
$ws.Range("R7").Value = 4724
$ws.Range("S7").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5274461', 'current_grade': None, 'current_score': 32.82, 'final_grade': None, 'final_score': 27.72, 'unposted_current_score': 32.82, 'unposted_current_grade': None, 'unposted_final_score': 27.72, 'unposted_final_grade': None}"

$ws.Range("P8").Value = "2025-07-20T13:56:14Z"
$ws.Range("R8").Value = 31390
$ws.Range("S8").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5046387', 'current_grade': None, 'current_score': 95.59, 'final_grade': None, 'final_score': 86.03, 'unposted_current_score': 95.59, 'unposted_current_grade': None, 'unposted_final_score': 86.03, 'unposted_final_grade': None}"

$ws.Range("P9").Value = "2025-07-21T00:31:41Z"
$ws.Range("R9").Value = 49690
$ws.Range("S9").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5018111', 'current_grade': None, 'current_score': 94.59, 'final_grade': None, 'final_score': 85.13, 'unposted_current_score': 94.59, 'unposted_current_grade': None, 'unposted_final_score': 85.13, 'unposted_final_grade': None}"

$ws.Range("P10").Value = "2025-07-19T15:43:54Z"
$ws.Range("R10").Value = 25964
$ws.Range("S10").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5341511', 'current_grade': None, 'current_score': 99.75, 'final_grade': None, 'final_score': 89.78, 'unposted_current_score': 99.75, 'unposted_current_grade': None, 'unposted_final_score': 89.78, 'unposted_final_grade': None}"

$ws.Range("P11").Value = "2025-07-18T19:27:30Z"
$ws.Range("R11").Value = 23578
$ws.Range("S11").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5254955', 'current_grade': None, 'current_score': 66.7, 'final_grade': None, 'final_score': 60.03, 'unposted_current_score': 66.7, 'unposted_current_grade': None, 'unposted_final_score': 60.03, 'unposted_final_grade': None}"

$ws.Range("P12").Value = "2025-07-21T14:04:14Z"
$ws.Range("R12").Value = 41358
$ws.Range("S12").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5078245', 'current_grade': None, 'current_score': 95.94, 'final_grade': None, 'final_score': 79.51, 'unposted_current_score': 95.94, 'unposted_current_grade': None, 'unposted_final_score': 79.51, 'unposted_final_grade': None}"

$ws.Range("P13").Value = "2025-07-21T08:37:59Z"
$ws.Range("R13").Value = 40803
$ws.Range("S13").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5226673', 'current_grade': None, 'current_score': 81.33, 'final_grade': None, 'final_score': 66.36, 'unposted_current_score': 81.33, 'unposted_current_grade': None, 'unposted_final_score': 66.36, 'unposted_final_grade': None}"

$ws.Range("P14").Value = "2025-07-20T17:32:25Z"
$ws.Range("R14").Value = 46093
$ws.Range("S14").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5285177', 'current_grade': None, 'current_score': 99.77, 'final_grade': None, 'final_score': 82.84, 'unposted_current_score': 99.77, 'unposted_current_grade': None, 'unposted_final_score': 82.84, 'unposted_final_grade': None}"

$ws.Range("P15").Value = "2025-07-21T03:34:32Z"
$ws.Range("R15").Value = 22114
$ws.Range("S15").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5324167', 'current_grade': None, 'current_score': 99.46, 'final_grade': None, 'final_score': 89.51, 'unposted_current_score': 99.46, 'unposted_current_grade': None, 'unposted_final_score': 89.51, 'unposted_final_grade': None}"

$ws.Range("P16").Value = "2025-07-20T21:32:41Z"
$ws.Range("R16").Value = 15392
$ws.Range("S16").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5266413', 'current_grade': None, 'current_score': 99.44, 'final_grade': None, 'final_score': 89.49, 'unposted_current_score': 99.44, 'unposted_current_grade': None, 'unposted_final_score': 89.49, 'unposted_final_grade': None}"

$ws.Range("P17").Value = "2025-07-20T17:22:55Z"
$ws.Range("R17").Value = 5927
$ws.Range("S17").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/4845467', 'current_grade': None, 'current_score': 53.82, 'final_grade': None, 'final_score': 48.43, 'unposted_current_score': 53.82, 'unposted_current_grade': None, 'unposted_final_score': 48.43, 'unposted_final_grade': None}"

$ws.Range("P18").Value = "2025-07-18T17:40:52Z"
$ws.Range("R18").Value = 13071
$ws.Range("S18").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/4596715', 'current_grade': None, 'current_score': 94.95, 'final_grade': None, 'final_score': 85.45, 'unposted_current_score': 94.95, 'unposted_current_grade': None, 'unposted_final_score': 85.45, 'unposted_final_grade': None}"

$ws.Range("P19").Value = "2025-07-21T00:31:58Z"
$ws.Range("R19").Value = 16215
$ws.Range("S19").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5285071', 'current_grade': None, 'current_score': 99.48, 'final_grade': None, 'final_score': 89.53, 'unposted_current_score': 99.48, 'unposted_current_grade': None, 'unposted_final_score': 89.53, 'unposted_final_grade': None}"

$ws.Range("P20").Value = "2025-07-19T23:17:52Z"
$ws.Range("R20").Value = 9803
$ws.Range("S20").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5314605', 'current_grade': None, 'current_score': 91.27, 'final_grade': None, 'final_score': 76.94, 'unposted_current_score': 91.27, 'unposted_current_grade': None, 'unposted_final_score': 76.94, 'unposted_final_grade': None}"

$ws.Range("P21").Value = "2025-07-18T07:12:22Z"
$ws.Range("R21").Value = 18501
$ws.Range("S21").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5341297', 'current_grade': None, 'current_score': 99.94, 'final_grade': None, 'final_score': 82.95, 'unposted_current_score': 99.94, 'unposted_current_grade': None, 'unposted_final_score': 82.95, 'unposted_final_grade': None}"

$ws.Range("P22").Value = "2025-07-18T16:10:06Z"
$ws.Range("R22").Value = 27150
$ws.Range("S22").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5322963', 'current_grade': None, 'current_score': 98.23, 'final_grade': None, 'final_score': 88.41, 'unposted_current_score': 98.23, 'unposted_current_grade': None, 'unposted_final_score': 88.41, 'unposted_final_grade': None}"

$ws.Range("P23").Value = "2025-07-21T04:16:25Z"
$ws.Range("R23").Value = 16257
$ws.Range("S23").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5234779', 'current_grade': None, 'current_score': 99.98, 'final_grade': None, 'final_score': 89.98, 'unposted_current_score': 99.98, 'unposted_current_grade': None, 'unposted_final_score': 89.98, 'unposted_final_grade': None}"

$ws.Range("P24").Value = "2025-07-16T19:10:03Z"
$ws.Range("R24").Value = 15700
$ws.Range("S24").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5244275', 'current_grade': None, 'current_score': 98.36, 'final_grade': None, 'final_score': 88.52, 'unposted_current_score': 98.36, 'unposted_current_grade': None, 'unposted_final_score': 88.52, 'unposted_final_grade': None}"

$ws.Range("P25").Value = "2025-07-19T00:21:36Z"
$ws.Range("R25").Value = 23516
$ws.Range("S25").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5033061', 'current_grade': None, 'current_score': 76.36, 'final_grade': None, 'final_score': 65.2, 'unposted_current_score': 76.36, 'unposted_current_grade': None, 'unposted_final_score': 65.2, 'unposted_final_grade': None}"

$ws.Range("P26").Value = "2025-07-21T06:31:06Z"
$ws.Range("S26").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5287447', 'current_grade': None, 'current_score': 99.17, 'final_grade': None, 'final_score': 89.25, 'unposted_current_score': 99.17, 'unposted_current_grade': None, 'unposted_final_score': 89.25, 'unposted_final_grade': None}"

$ws.Range("P27").Value = "2025-07-19T15:27:28Z"
$ws.Range("R27").Value = 12698
$ws.Range("S27").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5277805', 'current_grade': None, 'current_score': 99.34, 'final_grade': None, 'final_score': 82.45, 'unposted_current_score': 99.34, 'unposted_current_grade': None, 'unposted_final_score': 82.45, 'unposted_final_grade': None}"

$ws.Range("P28").Value = "2025-07-21T01:15:33Z"
$ws.Range("R28").Value = 21082
$ws.Range("S28").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5268281', 'current_grade': None, 'current_score': 90.25, 'final_grade': None, 'final_score': 81.22, 'unposted_current_score': 90.25, 'unposted_current_grade': None, 'unposted_final_score': 81.22, 'unposted_final_grade': None}"

$ws.Range("S29").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5321711', 'current_grade': None, 'current_score': 19.25, 'final_grade': None, 'final_score': 19.25, 'unposted_current_score': 19.25, 'unposted_current_grade': None, 'unposted_final_score': 19.25, 'unposted_final_grade': None}"

$ws.Range("P30").Value = "2025-07-16T22:03:51Z"
$ws.Range("R30").Value = 23966
$ws.Range("S30").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5342207', 'current_grade': None, 'current_score': 99.57, 'final_grade': None, 'final_score': 89.61, 'unposted_current_score': 99.57, 'unposted_current_grade': None, 'unposted_final_score': 89.61, 'unposted_final_grade': None}"

$ws.Range("P31").Value = "2025-07-21T00:41:17Z"
$ws.Range("R31").Value = 19095
$ws.Range("S31").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5255659', 'current_grade': None, 'current_score': 85.96, 'final_grade': None, 'final_score': 77.36, 'unposted_current_score': 85.96, 'unposted_current_grade': None, 'unposted_final_score': 77.36, 'unposted_final_grade': None}"

$ws.Range("P32").Value = "2025-07-21T00:12:40Z"
$ws.Range("S32").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5249109', 'current_grade': None, 'current_score': 60.1, 'final_grade': None, 'final_score': 55.48, 'unposted_current_score': 60.1, 'unposted_current_grade': None, 'unposted_final_score': 55.48, 'unposted_final_grade': None}"

$ws.Range("P33").Value = "2025-07-21T00:30:38Z"
$ws.Range("R33").Value = 26355
$ws.Range("S33").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/4995979', 'current_grade': None, 'current_score': 98.19, 'final_grade': None, 'final_score': 88.37, 'unposted_current_score': 98.19, 'unposted_current_grade': None, 'unposted_final_score': 88.37, 'unposted_final_grade': None}"

$ws.Range("P34").Value = "2025-07-19T21:25:01Z"
$ws.Range("R34").Value = 21667
$ws.Range("S34").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5260413', 'current_grade': None, 'current_score': 95.59, 'final_grade': None, 'final_score': 80.01, 'unposted_current_score': 95.59, 'unposted_current_grade': None, 'unposted_final_score': 80.01, 'unposted_final_grade': None}"

$ws.Range("P35").Value = "2025-07-21T00:47:15Z"
$ws.Range("S35").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5371577', 'current_grade': None, 'current_score': 99.76, 'final_grade': None, 'final_score': 89.78, 'unposted_current_score': 99.76, 'unposted_current_grade': None, 'unposted_final_score': 89.78, 'unposted_final_grade': None}"

$ws.Range("P36").Value = "2025-07-20T17:42:04Z"
$ws.Range("R36").Value = 11756
$ws.Range("S36").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5305513', 'current_grade': None, 'current_score': 76.52, 'final_grade': None, 'final_score': 64.74, 'unposted_current_score': 76.52, 'unposted_current_grade': None, 'unposted_final_score': 64.74, 'unposted_final_grade': None}"

$ws.Range("P37").Value = "2025-07-20T03:39:32Z"
$ws.Range("R37").Value = 21609
$ws.Range("S37").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5237497', 'current_grade': None, 'current_score': 91.37, 'final_grade': None, 'final_score': 75.91, 'unposted_current_score': 91.37, 'unposted_current_grade': None, 'unposted_final_score': 75.91, 'unposted_final_grade': None}"

$ws.Range("P38").Value = "2025-07-21T03:10:39Z"
$ws.Range("R38").Value = 50428
$ws.Range("S38").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5159167', 'current_grade': None, 'current_score': 70.59, 'final_grade': None, 'final_score': 59.92, 'unposted_current_score': 70.59, 'unposted_current_grade': None, 'unposted_final_score': 59.92, 'unposted_final_grade': None}"

$ws.Range("P39").Value = "2025-07-20T17:21:53Z"
$ws.Range("R39").Value = 19772
$ws.Range("S39").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5075195', 'current_grade': None, 'current_score': 95.3, 'final_grade': None, 'final_score': 78.77, 'unposted_current_score': 95.3, 'unposted_current_grade': None, 'unposted_final_score': 78.77, 'unposted_final_grade': None}"

$ws.Range("P40").Value = "2025-07-16T01:13:54Z"
$ws.Range("R40").Value = 30112
$ws.Range("S40").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5266969', 'current_grade': None, 'current_score': 94.67, 'final_grade': None, 'final_score': 85.21, 'unposted_current_score': 94.67, 'unposted_current_grade': None, 'unposted_final_score': 85.21, 'unposted_final_grade': None}"

$ws.Range("P41").Value = "2025-07-21T00:54:22Z"
$ws.Range("R41").Value = 29424
$ws.Range("S41").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5269555', 'current_grade': None, 'current_score': 99.15, 'final_grade': None, 'final_score': 89.24, 'unposted_current_score': 99.15, 'unposted_current_grade': None, 'unposted_final_score': 89.24, 'unposted_final_grade': None}"

$ws.Range("P42").Value = "2025-07-21T03:17:57Z"
$ws.Range("R42").Value = 14917
$ws.Range("S42").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5368843', 'current_grade': None, 'current_score': 84.65, 'final_grade': None, 'final_score': 71.96, 'unposted_current_score': 84.65, 'unposted_current_grade': None, 'unposted_final_score': 71.96, 'unposted_final_grade': None}"

$ws.Range("P43").Value = "2025-07-16T23:09:29Z"
$ws.Range("S43").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5264859', 'current_grade': None, 'current_score': 70.81, 'final_grade': None, 'final_score': 64.23, 'unposted_current_score': 70.81, 'unposted_current_grade': None, 'unposted_final_score': 64.23, 'unposted_final_grade': None}"

$ws.Range("P44").Value = "2025-07-21T05:23:59Z"
$ws.Range("R44").Value = 16019
$ws.Range("S44").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5081859', 'current_grade': None, 'current_score': 90.16, 'final_grade': None, 'final_score': 74.68, 'unposted_current_score': 90.16, 'unposted_current_grade': None, 'unposted_final_score': 74.68, 'unposted_final_grade': None}"

$ws.Range("P45").Value = "2025-07-18T17:24:12Z"
$ws.Range("R45").Value = 4894
$ws.Range("S45").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5291005', 'current_grade': None, 'current_score': 64.84, 'final_grade': None, 'final_score': 64.84, 'unposted_current_score': 64.84, 'unposted_current_grade': None, 'unposted_final_score': 64.84, 'unposted_final_grade': None}"

$ws.Range("P46").Value = "2025-07-21T13:22:48Z"
$ws.Range("R46").Value = 40917
$ws.Range("S46").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5245161', 'current_grade': None, 'current_score': 94.56, 'final_grade': None, 'final_score': 85.1, 'unposted_current_score': 94.56, 'unposted_current_grade': None, 'unposted_final_score': 85.1, 'unposted_final_grade': None}"

$ws.Range("P47").Value = "2025-07-19T06:38:16Z"
$ws.Range("R47").Value = 12504
$ws.Range("S47").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5264631', 'current_grade': None, 'current_score': 75.76, 'final_grade': None, 'final_score': 63.17, 'unposted_current_score': 75.76, 'unposted_current_grade': None, 'unposted_final_score': 63.17, 'unposted_final_grade': None}"

$ws.Range("P48").Value = "2025-07-18T16:34:01Z"
$ws.Range("R48").Value = 20599
$ws.Range("S48").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/4803449', 'current_grade': None, 'current_score': 86.03, 'final_grade': None, 'final_score': 71.85, 'unposted_current_score': 86.03, 'unposted_current_grade': None, 'unposted_final_score': 71.85, 'unposted_final_grade': None}"

$ws.Range("P49").Value = "2025-07-21T00:55:06Z"
$ws.Range("R49").Value = 25899
$ws.Range("S49").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5052305', 'current_grade': None, 'current_score': 86.29, 'final_grade': None, 'final_score': 72.76, 'unposted_current_score': 86.29, 'unposted_current_grade': None, 'unposted_final_score': 72.76, 'unposted_final_grade': None}"

$ws.Range("P50").Value = "2025-07-16T18:32:28Z"
$ws.Range("R50").Value = 31290
$ws.Range("S50").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5259461', 'current_grade': None, 'current_score': 97.04, 'final_grade': None, 'final_score': 87.34, 'unposted_current_score': 97.04, 'unposted_current_grade': None, 'unposted_final_score': 87.34, 'unposted_final_grade': None}"

$ws.Range("S51").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5371709', 'current_grade': None, 'current_score': 3.25, 'final_grade': None, 'final_score': 3.25, 'unposted_current_score': 3.25, 'unposted_current_grade': None, 'unposted_final_score': 3.25, 'unposted_final_grade': None}"

$ws.Range("P52").Value = "2025-07-20T02:50:42Z"
$ws.Range("R52").Value = 32179
$ws.Range("S52").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5333367', 'current_grade': None, 'current_score': 92.65, 'final_grade': None, 'final_score': 77.02, 'unposted_current_score': 92.65, 'unposted_current_grade': None, 'unposted_final_score': 77.02, 'unposted_final_grade': None}"

$ws.Range("P53").Value = "2025-07-18T14:52:56Z"
$ws.Range("S53").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5400867', 'current_grade': None, 'current_score': 83.78, 'final_grade': None, 'final_score': 83.78, 'unposted_current_score': 83.78, 'unposted_current_grade': None, 'unposted_final_score': 83.78, 'unposted_final_grade': None}"

$ws.Range("P54").Value = "2025-07-18T16:13:28Z"
$ws.Range("R54").Value = 14989
$ws.Range("S54").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5263813', 'current_grade': None, 'current_score': 88.28, 'final_grade': None, 'final_score': 79.45, 'unposted_current_score': 88.28, 'unposted_current_grade': None, 'unposted_final_score': 79.45, 'unposted_final_grade': None}"

$ws.Range("P55").Value = "2025-07-20T12:43:59Z"
$ws.Range("R55").Value = 13521
$ws.Range("S55").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5328305', 'current_grade': None, 'current_score': 98.77, 'final_grade': None, 'final_score': 88.9, 'unposted_current_score': 98.77, 'unposted_current_grade': None, 'unposted_final_score': 88.9, 'unposted_final_grade': None}"

$ws.Range("P56").Value = "2025-07-21T02:16:21Z"
$ws.Range("R56").Value = 28352
$ws.Range("S56").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5258093', 'current_grade': None, 'current_score': 97.71, 'final_grade': None, 'final_score': 87.94, 'unposted_current_score': 97.71, 'unposted_current_grade': None, 'unposted_final_score': 87.94, 'unposted_final_grade': None}"

$ws.Range("P57").Value = "2025-07-21T10:55:38Z"
$ws.Range("R57").Value = 76345
$ws.Range("S57").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5242323', 'current_grade': None, 'current_score': 76.79, 'final_grade': None, 'final_score': 76.79, 'unposted_current_score': 76.79, 'unposted_current_grade': None, 'unposted_final_score': 76.79, 'unposted_final_grade': None}"

$ws.Range("P58").Value = "2025-07-21T14:12:37Z"
$ws.Range("R58").Value = 22294
$ws.Range("S58").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5042553', 'current_grade': None, 'current_score': 94.71, 'final_grade': None, 'final_score': 78.29, 'unposted_current_score': 94.71, 'unposted_current_grade': None, 'unposted_final_score': 78.29, 'unposted_final_grade': None}"

$ws.Range("P59").Value = "2025-07-15T21:05:08Z"
$ws.Range("S59").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5385177', 'current_grade': None, 'current_score': 40.51, 'final_grade': None, 'final_score': 39.39, 'unposted_current_score': 40.51, 'unposted_current_grade': None, 'unposted_final_score': 39.39, 'unposted_final_grade': None}"

$ws.Range("P60").Value = "2025-07-19T13:20:18Z"
$ws.Range("R60").Value = 27579
$ws.Range("S60").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5343251', 'current_grade': None, 'current_score': 89.34, 'final_grade': None, 'final_score': 74.41, 'unposted_current_score': 89.34, 'unposted_current_grade': None, 'unposted_final_score': 74.41, 'unposted_final_grade': None}"

$ws.Range("P61").Value = "2025-07-18T17:35:46Z"
$ws.Range("R61").Value = 8052
$ws.Range("S61").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5277021', 'current_grade': None, 'current_score': 87.97, 'final_grade': None, 'final_score': 87.97, 'unposted_current_score': 87.97, 'unposted_current_grade': None, 'unposted_final_score': 87.97, 'unposted_final_grade': None}"

$ws.Range("P62").Value = "2025-07-20T23:18:32Z"
$ws.Range("R62").Value = 15565
$ws.Range("S62").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5418569', 'current_grade': None, 'current_score': 94.09, 'final_grade': None, 'final_score': 78.89, 'unposted_current_score': 94.09, 'unposted_current_grade': None, 'unposted_final_score': 78.89, 'unposted_final_grade': None}"

$ws.Range("P63").Value = "2025-07-21T12:37:24Z"
$ws.Range("R63").Value = 26642
$ws.Range("S63").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5237007', 'current_grade': None, 'current_score': 97.55, 'final_grade': None, 'final_score': 87.79, 'unposted_current_score': 97.55, 'unposted_current_grade': None, 'unposted_final_score': 87.79, 'unposted_final_grade': None}"

$ws.Range("P64").Value = "2025-07-21T00:30:14Z"
$ws.Range("S64").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5332285', 'current_grade': None, 'current_score': 98.01, 'final_grade': None, 'final_score': 88.21, 'unposted_current_score': 98.01, 'unposted_current_grade': None, 'unposted_final_score': 88.21, 'unposted_final_grade': None}"

$ws.Range("P65").Value = "2025-07-16T23:58:35Z"
$ws.Range("R65").Value = 21727
$ws.Range("S65").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5230621', 'current_grade': None, 'current_score': 99.31, 'final_grade': None, 'final_score': 89.38, 'unposted_current_score': 99.31, 'unposted_current_grade': None, 'unposted_final_score': 89.38, 'unposted_final_grade': None}"

$ws.Range("P66").Value = "2025-07-21T00:27:54Z"
$ws.Range("R66").Value = 1141
$ws.Range("S66").Value = "{'html_url': 'https://usflearn.instructure.com/courses/1962040/grades/5537663', 'current_grade': None, 'current_score': 42.86, 'final_grade': None, 'final_score': 15.0, 'unposted_current_score': 42.86, 'unposted_current_grade': None, 'unposted_final_score': 15.0, 'unposted_final_grade': None}"

$ws.Range("P67").Value = "2025-07-11T20:25:58Z"

$ws.Range("P68").Value = "2025-07-19T18:46:47Z"
$ws.Range("R68").Value = 7734

$ws.Range("P69").Value = "2025-07-21T15:10:21Z"
$ws.Range("R69").Value = 2203

$ws.Range("P70").Value = "2025-07-20T21:56:59Z"
$ws.Range("R70").Value = 43258

$ws.Range("P72").Value = "2025-07-21T00:28:01Z"
$ws.Range("R72").Value = 117150
